$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties labels in AD1:AF1, copying the
# existing bold/bordered/centered header formatting from AC1 (same cell
# style as the rest of row 1) and then setting the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows: every row from 2 to 46 gets the same team record (71-91-0
# wins/losses/ties) in the new AD/AE/AF columns.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
